# Insert a new data row at row 215 (pushing the existing rows 215-307 down
# to 216-308) and populate it with a new Sandia (watermelon) price record
# for "Vega Modelo de Temuco".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 215:307 down by one row, leaving a blank row 215 behind
# (Excel copies the formatting of the row above, matching style s="2" on D).
$ws.Rows(215).Insert()

$ws.Range("A215").Value = 10
$ws.Range("B215").Value = "Vega Modelo de Temuco"
$ws.Range("C215").Value = "La Araucanía"
$ws.Range("D215").Value = 44452
$ws.Range("E215").Value = 9
$ws.Range("F215").Value = 100112028
$ws.Range("G215").Value = "Sandia"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 200
$ws.Range("K215").Value = 3200
$ws.Range("L215").Value = 3200
$ws.Range("M215").Value = 3200
$ws.Range("N215").Value = "`$/unidad"
$ws.Range("O215").Value = "Brasil"
$ws.Range("P215").Value = 3200
$ws.Range("Q215").Value = 1
$ws.Range("R215").Value = "Hortaliza"
